# Add a new "Sources" slide at the end of the deck (Title and Content layout,
# same as every other content slide in this deck).
$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# --- Title placeholder -----------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Sources"
$title.LanguageID = "en-GB"

# --- Body / content placeholder --------------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange

$datasetLink = "https://www.kaggle.com/datasets/pavansubhasht/ibm-hr-analytics-attrition-dataset/discussion/233758"
$helpLink = "https://inseaddataanalytics.github.io/INSEADAnalytics/groupprojects/January2018FBL/IBM_Attrition_VSS.html"

# Paragraph 1: "Dataset: " followed by a hyperlinked URL run.
$body.Text = "Dataset: "
$body.LanguageID = "en-GB"
[void]$body.InsertAfter($datasetLink)
$body.LanguageID = "en-GB"

# Paragraph 2: "Help" + ": " + a hyperlinked URL run.
[void]$body.InsertAfter([char]13 + "Help")
$body.LanguageID = "en-GB"
[void]$body.InsertAfter(": ")
$body.LanguageID = "en-GB"
[void]$body.InsertAfter($helpLink)
$body.LanguageID = "en-GB"

# Paragraph 3: empty trailing paragraph.
[void]$body.InsertAfter([char]13)
$body.LanguageID = "en-GB"

# --- Turn the two URLs into real hyperlinks --------------------------------
$prefix1 = "Dataset: "
$urlRange1 = $body.Characters($prefix1.Length + 1, $datasetLink.Length)
$urlRange1.ActionSettings.Item(1).Hyperlink.Address = $datasetLink

$line2Offset = $prefix1.Length + $datasetLink.Length + 1
$prefix2 = "Help: "
$urlRange2 = $body.Characters($line2Offset + $prefix2.Length + 1, $helpLink.Length)
$urlRange2.ActionSettings.Item(1).Hyperlink.Address = $helpLink
